$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 4
$ws.Cells.Item(4, 8).Value = 1000
$ws.Cells.Item(4, 9).Value = 1000
$ws.Cells.Item(4, 11).Value = 1000
$ws.Cells.Item(4, 13).Value = -886
# ALC row 8
$ws.Cells.Item(8, 8).Value = 203.66667
$ws.Cells.Item(8, 9).Value = 203.66667
$ws.Cells.Item(8, 10).Value = 0
$ws.Cells.Item(8, 11).Value = 611.00001
$ws.Cells.Item(8, 12).Value = 0
$ws.Cells.Item(8, 13).Value = -472.00001
$ws.Cells.Item(8, 14).Value = ""
# ALC row 31
$ws.Cells.Item(31, 8).Value = 878.3333
$ws.Cells.Item(31, 9).Value = 841.5
$ws.Cells.Item(31, 11).Value = 2524.5
$ws.Cells.Item(31, 13).Value = -2294.5
# ALC row 33
$ws.Cells.Item(33, 8).Value = 226.85715
$ws.Cells.Item(33, 9).Value = 226.85715
$ws.Cells.Item(33, 11).Value = 226.85715
$ws.Cells.Item(33, 13).Value = 2.14285000000001
# ALC row 132
$ws.Cells.Item(132, 8).Value = 239339.92
$ws.Cells.Item(132, 9).Value = 1086.3158
$ws.Cells.Item(132, 11).Value = 3258.9474
$ws.Cells.Item(132, 13).Value = -728.9474

$ws = $wb.Worksheets.Item("ARM")
# ARM row 2
$ws.Cells.Item(2, 8).Value = 1302.8975
$ws.Cells.Item(2, 9).Value = 1310.8948
$ws.Cells.Item(2, 11).Value = 1310.8948
$ws.Cells.Item(2, 13).Value = -1197.8948
# ARM row 74
$ws.Cells.Item(74, 8).Value = 595
$ws.Cells.Item(74, 9).Value = 595
$ws.Cells.Item(74, 11).Value = 595
$ws.Cells.Item(74, 13).Value = 279
# ARM row 77
$ws.Cells.Item(77, 8).Value = 595
$ws.Cells.Item(77, 9).Value = 595
$ws.Cells.Item(77, 11).Value = 2975
$ws.Cells.Item(77, 13).Value = 1393
# ARM row 116
$ws.Cells.Item(116, 8).Value = 1302.8975
$ws.Cells.Item(116, 9).Value = 1310.8948
$ws.Cells.Item(116, 11).Value = 1310.8948
$ws.Cells.Item(116, 13).Value = 983.1052

$ws = $wb.Worksheets.Item("BSM")
# BSM row 3
$ws.Cells.Item(3, 8).Value = 1302.8975
$ws.Cells.Item(3, 9).Value = 1310.8948
$ws.Cells.Item(3, 11).Value = 1310.8948
$ws.Cells.Item(3, 13).Value = -1196.8948

$ws = $wb.Worksheets.Item("CRP")
# CRP row 86
$ws.Cells.Item(86, 8).Value = 3616
$ws.Cells.Item(86, 9).Value = 3538.7144
$ws.Cells.Item(86, 11).Value = 3538.7144
$ws.Cells.Item(86, 13).Value = -2415.7144
# CRP row 89
$ws.Cells.Item(89, 8).Value = 3616
$ws.Cells.Item(89, 9).Value = 3538.7144
$ws.Cells.Item(89, 11).Value = 17693.572
$ws.Cells.Item(89, 13).Value = -12077.572
# CRP row 99
$ws.Cells.Item(99, 8).Value = 1431.8667
$ws.Cells.Item(99, 9).Value = 1312.3
$ws.Cells.Item(99, 10).Value = 1671
$ws.Cells.Item(99, 11).Value = 1312.3
$ws.Cells.Item(99, 12).Value = 1671
$ws.Cells.Item(99, 13).Value = 185.7
$ws.Cells.Item(99, 14).Value = -4667
# CRP row 126
$ws.Cells.Item(126, 8).Value = 1431.8667
$ws.Cells.Item(126, 9).Value = 1312.3
$ws.Cells.Item(126, 10).Value = 1671
$ws.Cells.Item(126, 11).Value = 3936.9
$ws.Cells.Item(126, 12).Value = 5013
$ws.Cells.Item(126, 13).Value = -1466.9
$ws.Cells.Item(126, 14).Value = -9953

$ws = $wb.Worksheets.Item("CUL")
# CUL row 17
$ws.Cells.Item(17, 8).Value = 3465.75
$ws.Cells.Item(17, 9).Value = 625
$ws.Cells.Item(17, 11).Value = 1875
$ws.Cells.Item(17, 13).Value = -1706
# CUL row 60
$ws.Cells.Item(60, 8).Value = 497.77777
$ws.Cells.Item(60, 9).Value = 497.77777
$ws.Cells.Item(60, 11).Value = 1493.33331
$ws.Cells.Item(60, 13).Value = -1242.33331
# CUL row 69
$ws.Cells.Item(69, 8).Value = 10999
$ws.Cells.Item(69, 9).Value = 9998
$ws.Cells.Item(69, 10).Value = 12000
$ws.Cells.Item(69, 11).Value = 29994
$ws.Cells.Item(69, 12).Value = 36000
$ws.Cells.Item(69, 13).Value = -29183
$ws.Cells.Item(69, 14).Value = -37622
# CUL row 72
$ws.Cells.Item(72, 8).Value = 10999
$ws.Cells.Item(72, 9).Value = 9998
$ws.Cells.Item(72, 10).Value = 12000
$ws.Cells.Item(72, 11).Value = 89982
$ws.Cells.Item(72, 12).Value = 108000
$ws.Cells.Item(72, 13).Value = -85926
$ws.Cells.Item(72, 14).Value = -116112
# CUL row 112
$ws.Cells.Item(112, 8).Value = 0
$ws.Cells.Item(112, 9).Value = 0
$ws.Cells.Item(112, 11).Value = 0
$ws.Cells.Item(112, 13).Value = ""
# CUL row 121
$ws.Cells.Item(121, 8).Value = 85170.75
$ws.Cells.Item(121, 10).Value = 2113.2222
$ws.Cells.Item(121, 12).Value = 6339.6666
$ws.Cells.Item(121, 14).Value = -8959.6666
# CUL row 131
$ws.Cells.Item(131, 8).Value = 18958.223
$ws.Cells.Item(131, 9).Value = 56169.3
$ws.Cells.Item(131, 10).Value = 1650.7441
$ws.Cells.Item(131, 11).Value = 168507.9
$ws.Cells.Item(131, 12).Value = 4952.2323
$ws.Cells.Item(131, 13).Value = -163467.9
$ws.Cells.Item(131, 14).Value = -15032.2323

$ws = $wb.Worksheets.Item("LTW")
# LTW row 22
$ws.Cells.Item(22, 8).Value = 1647.75
$ws.Cells.Item(22, 9).Value = 1424.5
$ws.Cells.Item(22, 10).Value = 2094.25
$ws.Cells.Item(22, 11).Value = 1424.5
$ws.Cells.Item(22, 12).Value = 2094.25
$ws.Cells.Item(22, 13).Value = -1129.5
$ws.Cells.Item(22, 14).Value = -2684.25
# LTW row 27
$ws.Cells.Item(27, 8).Value = 1647.75
$ws.Cells.Item(27, 9).Value = 1424.5
$ws.Cells.Item(27, 10).Value = 2094.25
$ws.Cells.Item(27, 11).Value = 1424.5
$ws.Cells.Item(27, 12).Value = 2094.25
$ws.Cells.Item(27, 13).Value = -1317.5
$ws.Cells.Item(27, 14).Value = -2308.25
# LTW row 40
$ws.Cells.Item(40, 8).Value = 4214
$ws.Cells.Item(40, 9).Value = 4115.875
$ws.Cells.Item(40, 11).Value = 4115.875
$ws.Cells.Item(40, 13).Value = -3979.875
# LTW row 46
$ws.Cells.Item(46, 8).Value = 1280.7916
$ws.Cells.Item(46, 9).Value = 1203.2142
$ws.Cells.Item(46, 10).Value = 1389.4
$ws.Cells.Item(46, 11).Value = 1203.2142
$ws.Cells.Item(46, 12).Value = 1389.4
$ws.Cells.Item(46, 13).Value = -1015.2142
$ws.Cells.Item(46, 14).Value = -1765.4
# LTW row 61
$ws.Cells.Item(61, 8).Value = 1792.3125
$ws.Cells.Item(61, 9).Value = 1792.3125
$ws.Cells.Item(61, 10).Value = 0
$ws.Cells.Item(61, 11).Value = 1792.3125
$ws.Cells.Item(61, 12).Value = 0
$ws.Cells.Item(61, 13).Value = -1590.3125
$ws.Cells.Item(61, 14).Value = ""
# LTW row 100
$ws.Cells.Item(100, 8).Value = 3522.8333
$ws.Cells.Item(100, 9).Value = 2050
$ws.Cells.Item(100, 11).Value = 2050
$ws.Cells.Item(100, 13).Value = -1509
# LTW row 113
$ws.Cells.Item(113, 8).Value = 1792.3125
$ws.Cells.Item(113, 9).Value = 1792.3125
$ws.Cells.Item(113, 10).Value = 0
$ws.Cells.Item(113, 11).Value = 1792.3125
$ws.Cells.Item(113, 12).Value = 0
$ws.Cells.Item(113, 13).Value = 377.6875
$ws.Cells.Item(113, 14).Value = ""
# LTW row 122
$ws.Cells.Item(122, 8).Value = 3046.1052
$ws.Cells.Item(122, 9).Value = 2837.4814
$ws.Cells.Item(122, 11).Value = 8512.4442
$ws.Cells.Item(122, 13).Value = -6062.4442
# LTW row 125
$ws.Cells.Item(125, 8).Value = 0
$ws.Cells.Item(125, 10).Value = 0
$ws.Cells.Item(125, 12).Value = 0
$ws.Cells.Item(125, 14).Value = ""
# LTW row 136
$ws.Cells.Item(136, 8).Value = 3998.4285
$ws.Cells.Item(136, 9).Value = 3998.1667
$ws.Cells.Item(136, 11).Value = 11994.5001
$ws.Cells.Item(136, 13).Value = -9444.500100000001

$ws = $wb.Worksheets.Item("WVR")
# WVR row 93
$ws.Cells.Item(93, 8).Value = 0
$ws.Cells.Item(93, 10).Value = 0
$ws.Cells.Item(93, 12).Value = 0
$ws.Cells.Item(93, 14).Value = ""

